# "After teamwork defense commit. Minor fixes."
# Correct the "Needed Kits" / "Needed Boots" counts and the resulting
# "Total Costs" for every team on the Teams Info sheet. These columns are
# stored as text in the workbook, so a leading apostrophe is used to keep
# the entry literal (matches the existing "number-looking text" cells)
# instead of letting Excel auto-convert it to a numeric value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# CSKA
$ws.Range("H2").Value = "'60"
$ws.Range("I2").Value = "'60"
$ws.Range("L2").Value = "'28200"

# Levski
$ws.Range("H3").Value = "'44"
$ws.Range("I3").Value = "'44"
$ws.Range("L3").Value = "'22000"

# Beroe
$ws.Range("H4").Value = "'60"
$ws.Range("I4").Value = "'60"
$ws.Range("L4").Value = "'18000"

# Ludogorets
$ws.Range("H5").Value = "'44"
$ws.Range("I5").Value = "'44"
$ws.Range("L5").Value = "'21120"

# Litex
$ws.Range("H6").Value = "'44"
$ws.Range("I6").Value = "'44"
$ws.Range("L6").Value = "'16720"
